# Update the "as_of_utc" timestamp column (AA) for rows 2-26 on the
# "Главные" and "Линейные" sheets from 2025-11-22 03:03:39 to
# 2025-11-22 07:04:51.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Главные", "Линейные")
$newTimestamp = "2025-11-22 07:04:51"

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Range("AA$row").Value = $newTimestamp
    }
}
